# Update the "想去人数" (want-to-go count) figures in column F on the
# "展览", "演出" and "全部类型" sheets to the newly regenerated values.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsShow    = $wb.Worksheets.Item("演出")
$wsAll     = $wb.Worksheets.Item("全部类型")

# 展览 ("Exhibitions") sheet
$wsExhibit.Range("F3").Value  = 384
$wsExhibit.Range("F5").Value  = 1305
$wsExhibit.Range("F7").Value  = 2482
$wsExhibit.Range("F8").Value  = 889
$wsExhibit.Range("F9").Value  = 18597
$wsExhibit.Range("F11").Value = 1902
$wsExhibit.Range("F15").Value = 601
$wsExhibit.Range("F18").Value = 70

# 演出 ("Shows") sheet
$wsShow.Range("F5").Value = 166
$wsShow.Range("F8").Value = 109

# 全部类型 ("All types") sheet - combined view of all rows
$wsAll.Range("F6").Value  = 384
$wsAll.Range("F10").Value = 1305
$wsAll.Range("F13").Value = 166
$wsAll.Range("F15").Value = 2482
$wsAll.Range("F16").Value = 889
$wsAll.Range("F17").Value = 18597
$wsAll.Range("F20").Value = 109
$wsAll.Range("F23").Value = 1902
$wsAll.Range("F28").Value = 601
$wsAll.Range("F32").Value = 70
